$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "Address Block" entry becomes "Data Block" entry (var[] / Data Block / Block of data) ---
$ws.Range("A8").Value = "var[]"
$ws.Range("B8").Value = "Data Block"
$ws.Range("C8").Value = "Block of data"
$ws.Range("D8").ClearContents()
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = "Allocated"

# --- Row 9 (old "Block of data" row) and row 10 (stray J10 value) removed ---
$ws.Range("A9").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("J10").ClearContents()

# --- New "Data Block" detail table starting at row 11 ---
$ws.Range("A11").Value = "Data Block"

$ws.Range("A12").Value = "64b"
$ws.Range("B12").Value = "u64"
$ws.Range("C12").Value = "total size of data"
$ws.Range("D12").Value = "size of data"

$ws.Range("A13").Value = "32b"
$ws.Range("B13").Value = "u32"
$ws.Range("C13").Value = "state flag"
$ws.Range("D13").Value = "if free we can allocate"
$ws.Range("E13").Value = "deletes can nullify this data"

$ws.Range("A14").Value = "64b"
$ws.Range("B14").Value = "u64"
$ws.Range("C14").Value = "address of next chunk"
$ws.Range("D14").ClearContents()
$ws.Range("E14").ClearContents()

$ws.Range("A15").Value = "32b"
$ws.Range("B15").Value = "checksum"
$ws.Range("C15").Value = "checksum, 32b probably sufficient 0 if not used"

$ws.Range("A16").Value = "var"
$ws.Range("B16").Value = "u8[]"
$ws.Range("C16").Value = "the data"

# --- Old rows 17 and 21-24 (size-of-chunk + Address Block sub-table) removed entirely ---
$ws.Range("A17").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()

$ws.Range("A21").ClearContents()
$ws.Range("A22").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("B23").ClearContents()
$ws.Range("A24").ClearContents()
$ws.Range("B24").ClearContents()

# --- Selection / active cell now sits on C16 (last populated cell) ---
$ws.Range("C16").Select()
